# "Robótica" -> "Compiladores" deck rewrite, per the commit's XML diff.
$p = $ppt.ActivePresentation

# --- Slide 1 (title slide): update title + author -------------------------
$s1 = $p.Slides.Item(1)
$s1.Shapes.Item(1).TextFrame.TextRange.Text = "Compiladores"
$s1.Shapes.Item(2).TextFrame.TextRange.Text = "Autor: Kézia Vasconcelos"

# --- Slide 2: Introdução ---------------------------------------------------
$s2 = $p.Slides.Item(2)
$s2.Shapes.Item(1).TextFrame.TextRange.Text = "Introdução"
$s2.Shapes.Item(2).TextFrame.TextRange.Text = "Definição de compiladores`nImportância na programação`nProcesso de compilação"

# --- Slide 3: Fases de um compilador ---------------------------------------
$s3 = $p.Slides.Item(3)
$s3.Shapes.Item(1).TextFrame.TextRange.Text = "Fases de um compilador"
$s3.Shapes.Item(2).TextFrame.TextRange.Text = "Análise léxica`nAnálise sintática`nAnálise semântica`nGeração de código intermediário`nOtimização de código`nGeração de código final"

# --- Slide 4: Análise léxica -------------------------------------------------
$s4 = $p.Slides.Item(4)
$s4.Shapes.Item(1).TextFrame.TextRange.Text = "Análise léxica"
$s4.Shapes.Item(2).TextFrame.TextRange.Text = "Identificação de tokens`nRemoção de espaços em branco e comentários`nGeração do código interno"

# --- Slide 5: Análise sintática ----------------------------------------------
$s5 = $p.Slides.Item(5)
$s5.Shapes.Item(1).TextFrame.TextRange.Text = "Análise sintática"
$s5.Shapes.Item(2).TextFrame.TextRange.Text = "Verificação da estrutura gramatical`nConstrução da árvore sintática`nDetecção de erros sintáticos"

# --- Slide 6: Análise semântica ----------------------------------------------
$s6 = $p.Slides.Item(6)
$s6.Shapes.Item(1).TextFrame.TextRange.Text = "Análise semântica"
$s6.Shapes.Item(2).TextFrame.TextRange.Text = "Verificação de tipos de dados`nResolução de ambiguidades`nVerificação de escopo de variáveis"

# --- Slide 7: Geração de código intermediário --------------------------------
$s7 = $p.Slides.Item(7)
$s7.Shapes.Item(1).TextFrame.TextRange.Text = "Geração de código intermediário"
$s7.Shapes.Item(2).TextFrame.TextRange.Text = "Representação intermédia`nFacilita a otimização do código`nPonte entre as fases de análise e de geração de código final"

# --- Slide 8: Otimização de código -------------------------------------------
$s8 = $p.Slides.Item(8)
$s8.Shapes.Item(1).TextFrame.TextRange.Text = "Otimização de código"
$s8.Shapes.Item(2).TextFrame.TextRange.Text = "Melhoria da eficiência do código`nRedução de redundâncias`nExemplo de otimizações comuns"

# --- Slide 9: Geração de código final -----------------------------------------
$s9 = $p.Slides.Item(9)
$s9.Shapes.Item(1).TextFrame.TextRange.Text = "Geração de código final"
$s9.Shapes.Item(2).TextFrame.TextRange.Text = "Tradução do código intermediário para a linguagem de máquina`nUso de registradores e instruções específicas`nProdução do executável final"

# --- Slide 10 (was "Conclusão" / robótica wrap-up): becomes the new
#     "Conclusão" slide with compiler-themed bullets. It keeps its identity
#     but four new slides get inserted in front of it below, which pushes it
#     down to its final position (14).
$s10 = $p.Slides.Item(10)
$s10.Shapes.Item(1).TextFrame.TextRange.Text = "Conclusão"
$s10.Shapes.Item(2).TextFrame.TextRange.Text = "Compiladores desempenham papel fundamental na programação`nProcesso de compilação é complexo e exige várias etapas`nImportância de conhecer o funcionamento dos compiladores para programadores"

# --- New slide 10: Exemplos de compiladores ---------------------------------
$sEx = $p.Slides.Add(10, 2)
$sEx.Shapes.Item(1).TextFrame.TextRange.Text = "Exemplos de compiladores"
$sEx.Shapes.Item(2).TextFrame.TextRange.Text = "GCC (GNU Compiler Collection)`nClang`nVisual Studio Compiler"

# --- New slide 11: Ferramentas auxiliares -----------------------------------
$sFer = $p.Slides.Add(11, 2)
$sFer.Shapes.Item(1).TextFrame.TextRange.Text = "Ferramentas auxiliares"
$sFer.Shapes.Item(2).TextFrame.TextRange.Text = "IDEs (Ambientes de Desenvolvimento Integrado)`nDepuradores (debuggers)`nPerfis de desempenho (profilers)"

# --- New slide 12: Desafios na construção de compiladores -------------------
$sDes = $p.Slides.Add(12, 2)
$sDes.Shapes.Item(1).TextFrame.TextRange.Text = "Desafios na construção de compiladores"
$sDes.Shapes.Item(2).TextFrame.TextRange.Text = "Lidar com ambiguidades na linguagem`nOtimização de código eficiente`nSuporte a múltiplas plataformas"

# --- New slide 13: Futuro dos compiladores ----------------------------------
$sFut = $p.Slides.Add(13, 2)
$sFut.Shapes.Item(1).TextFrame.TextRange.Text = "Futuro dos compiladores"
$sFut.Shapes.Item(2).TextFrame.TextRange.Text = "Avanços em otimização de código`nAumento da integração com IDEs`nSuporte a novas linguagens de programação"

# (slide 14 is now the old slide 10, relabeled "Conclusão" above, pushed down
# to index 14 by the four inserts.)

# --- New slide 15: Referências ----------------------------------------------
$sRef = $p.Slides.Add(15, 2)
$sRef.Shapes.Item(1).TextFrame.TextRange.Text = "Referências"
$sRef.Shapes.Item(2).TextFrame.TextRange.Text = "Livros sobre compiladores`nArtigos acadêmicos`nSites especializados em programação e compiladores"

# --- New slide 16: Perguntas -------------------------------------------------
$sPer = $p.Slides.Add(16, 2)
$sPer.Shapes.Item(1).TextFrame.TextRange.Text = "Perguntas"
$sPer.Shapes.Item(2).TextFrame.TextRange.Text = "Momento para esclarecer dúvidas`nDiscussão sobre o tema`nAgradecimentos."
